$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'5"
$ws.Range("D2").Value = 0.0883
$ws.Range("E2").Value = 0.09619999999999999
$ws.Range("F2").Value = 0.04600000000000001
$ws.Range("G2").Value = 0.1268072289156627
$ws.Range("H2").Value = 0.1268072289156627
$ws.Range("I2").Value = 0.1994164156626506
$ws.Range("J2").Value = 0.1570481353462966
$ws.Range("K2").Value = 157.89
$ws.Range("L2").Value = 0.2972326807228915
$ws.Range("M2").Value = 119.61
$ws.Range("N2").Value = 0.04232933432423824
$ws.Range("O2").Value = 0.7575527265817975
$ws.Range("P2").Value = 86.81
$ws.Range("Q2").Value = 0.0307215911101674
$ws.Range("R2").Value = 0.5498131610614986
$ws.Range("S2").Value = 32.8
$ws.Range("T2").Value = 0.274224563163615
$ws.Range("U2").Value = 190
$ws.Range("V2").Value = 0.06723997593516651
$ws.Range("W2").Value = 0.1149786660678195
$ws.Range("X2").Value = 0.03007806254876341
$ws.Range("Y2").Value = 0.08490060351905604
$ws.Range("Z2").Value = 0.1676857659659767
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.02560421492435847
$ws.Range("AC2").Value = -0.02560421492435847
$ws.Range("AD2").Value = 2526.01
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2526.01
$ws.Range("AG2").Value = 2336.01
$ws.Range("AH2").Value = 0.4720005381457516
$ws.Range("AI2").Value = 0.7364868607882069
$ws.Range("AJ2").Value = 0.4525651382971922
$ws.Range("AK2").Value = 0.721033023541504
$ws.Range("AL2").Value = 7.638999999999999
$ws.Range("AM2").Value = 7.638999999999999
$ws.Range("AN2").Value = 22.99508420573509
$ws.Range("AO2").Value = 13.86699829820657
$ws.Range("AP2").Value = 21.26545289030496
$ws.Range("AQ2").Value = 13.86699829820657
$ws.Range("D3").Value = 0.06660000000000001
$ws.Range("E3").Value = 0.09619999999999999
$ws.Range("F3").Value = 0.07200000000000001
$ws.Range("G3").Value = 0.3962940824865511
$ws.Range("H3").Value = 0.3962940824865511
$ws.Range("I3").Value = 0.6264196054991034
$ws.Range("J3").Value = 0.463163530706866
$ws.Range("K3").Value = 76.59999999999999
$ws.Range("L3").Value = 0.4578601315002988
$ws.Range("M3").Value = 53.3
$ws.Range("N3").Value = 0.03189897659943743
$ws.Range("O3").Value = 0.695822454308094
$ws.Range("P3").Value = 53.3
$ws.Range("Q3").Value = 0.03189897659943743
$ws.Range("R3").Value = 0.695822454308094
$ws.Range("U3").Value = 126.6
$ws.Range("V3").Value = 0.0757675504219283
$ws.Range("W3").Value = 0.4372146118721461
$ws.Range("X3").Value = 0.01896771930379554
$ws.Range("Y3").Value = 0.4182468925683506
$ws.Range("Z3").Value = 1.972877358490566
$ws.Range("AA3").Value = 0.9137648430101261
$ws.Range("AB3").Value = 0.0189721443182694
$ws.Range("AC3").Value = 0.8947926986918566
$ws.Range("AD3").Value = 2.01
$ws.Range("AF3").Value = 2.01
$ws.Range("AG3").Value = -124.59
$ws.Range("AH3").Value = 0.001201499184056524
$ws.Range("AI3").Value = 0.01047911996246285
$ws.Range("AJ3").Value = -0.08057245959736403
$ws.Range("AK3").Value = -1.910596534273884
$ws.Range("AL3").Value = 0.129
$ws.Range("AM3").Value = 0.129
$ws.Range("AN3").Value = 0.01866295264623955
$ws.Range("AO3").Value = 812.4031007751937
$ws.Range("AP3").Value = -1.156824512534819
$ws.Range("AQ3").Value = 812.4031007751937
$ws.Range("D4").Value = 0.197
$ws.Range("E4").Value = 0.0815
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6.75
$ws.Range("L4").Value = 0.2755102040816326
$ws.Range("M4").Value = 5.18
$ws.Range("N4").Value = 0.04832089552238806
$ws.Range("O4").Value = 0.7674074074074073
$ws.Range("P4").Value = 5.18
$ws.Range("Q4").Value = 0.04832089552238806
$ws.Range("R4").Value = 0.7674074074074073
$ws.Range("U4").Value = 11.5
$ws.Range("V4").Value = 0.1072761194029851
$ws.Range("W4").Value = 0.06736526946107785
$ws.Range("X4").Value = 0.02269774371605136
$ws.Range("Y4").Value = 0.04466752574502649
$ws.Range("Z4").Value = 0.1618978391594528
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.02268184509868598
$ws.Range("AC4").Value = -0.02268184509868598
$ws.Range("AD4").Value = 54.6
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 54.6
$ws.Range("AG4").Value = 43.1
$ws.Range("AH4").Value = 0.3374536464771323
$ws.Range("AI4").Value = 0.3475493316359007
$ws.Range("AJ4").Value = 0.2867598137059215
$ws.Range("AK4").Value = 0.2960164835164835
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("B5").Value = "RCE Capital Berhad (KLSE:RCECAP)"
$ws.Range("D5").Value = 0.11
$ws.Range("E5").Value = 0.216
$ws.Range("F5").Value = 0.02
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 27.7
$ws.Range("L5").Value = 0.5782881002087683
$ws.Range("M5").Value = 10.23
$ws.Range("N5").Value = 0.04182338511856092
$ws.Range("O5").Value = 0.3693140794223827
$ws.Range("P5").Value = 9.33
$ws.Range("Q5").Value = 0.03814390842191333
$ws.Range("R5").Value = 0.3368231046931408
$ws.Range("S5").Value = 0.9000000000000004
$ws.Range("T5").Value = 0.08797653958944285
$ws.Range("U5").Value = 12.6
$ws.Range("V5").Value = 0.05151267375306623
$ws.Range("W5").Value = 0.1842980705256154
$ws.Range("X5").Value = 0.03007806254876341
$ws.Range("Y5").Value = 0.154220007976852
$ws.Range("Z5").Value = 0.1008633396504527
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.02560421492435847
$ws.Range("AC5").Value = -0.02560421492435847
$ws.Range("AD5").Value = 370.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 370.5
$ws.Range("AG5").Value = 357.9
$ws.Range("AH5").Value = 0.6023410827507722
$ws.Range("AI5").Value = 0.6821948075860799
$ws.Range("AJ5").Value = 0.5940248962655601
$ws.Range("AK5").Value = 0.6746465598491989
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("B6").Value = "Johan Holdings Berhad (KLSE:JOHAN)"
$ws.Range("D6").Value = -0.171
$ws.Range("E6").ClearContents()
$ws.Range("G6").Value = 0.05120772946859904
$ws.Range("H6").Value = 0.05120772946859904
$ws.Range("I6").Value = 0.05458937198067632
$ws.Range("J6").Value = 0.05458937198067632
$ws.Range("K6").Value = -4.36
$ws.Range("L6").Value = -0.2106280193236715
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").ClearContents()
$ws.Range("U6").Value = 23.4
$ws.Range("V6").Value = 0.5749385749385748
$ws.Range("W6").Value = -0.2812903225806452
$ws.Range("X6").Value = 0.03964647239213256
$ws.Range("Y6").Value = -0.3209367949727777
$ws.Range("Z6").Value = 0.1375415282392027
$ws.Range("AA6").Value = 0.007508305647840531
$ws.Range("AB6").Value = 0.03314898889216905
$ws.Range("AC6").Value = -0.02564068324432852
$ws.Range("AD6").Value = 114.7
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 114.7
$ws.Range("AG6").Value = 91.30000000000001
$ws.Range("AH6").Value = 0.7380952380952381
$ws.Range("AI6").Value = 0.9095955590800952
$ws.Range("AJ6").Value = 0.6916666666666668
$ws.Range("AK6").Value = 0.8889970788704965
$ws.Range("AL6").Value = 7.51
$ws.Range("AM6").Value = 7.51
$ws.Range("AN6").Value = 53.34883720930233
$ws.Range("AO6").Value = 0.1504660452729694
$ws.Range("AP6").Value = 42.46511627906978
$ws.Range("AQ6").Value = 0.1504660452729694
$ws.Range("B7").Value = "AEON Credit Service (M) Berhad (KLSE:AEONCR)"
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("K7").Value = 51.2
$ws.Range("L7").Value = 0.189069423929099
$ws.Range("M7").Value = 50.9
$ws.Range("N7").Value = 0.06677161222615768
$ws.Range("O7").Value = 0.9941406249999999
$ws.Range("P7").Value = 19
$ws.Range("Q7").Value = 0.02492457037911584
$ws.Range("R7").Value = 0.37109375
$ws.Range("S7").Value = 31.9
$ws.Range("T7").Value = 0.6267190569744597
$ws.Range("U7").Value = 15.9
$ws.Range("V7").Value = 0.02085792994883904
$ws.Range("W7").Value = 0.1149786660678195
$ws.Range("X7").Value = 0.03806623105228851
$ws.Range("Y7").Value = 0.07691243501553094
$ws.Range("Z7").Value = 0.1174175085635
$ws.Range("AB7").Value = 0.02692928033325862
$ws.Range("AC7").Value = -0.02692928033325862
$ws.Range("AD7").Value = 1984.2
$ws.Range("AF7").Value = 1984.2
$ws.Range("AG7").Value = 1968.3
$ws.Range("AH7").Value = 0.7224467504096123
$ws.Range("AI7").Value = 0.8227391466600324
$ws.Range("AJ7").Value = 0.7208305866842453
$ws.Range("AK7").Value = 0.8215627347858752

$ws.Rows("8:8").Delete()
